$d = $word.ActiveDocument

# Remove the placeholder text "vnpt.SiteAddress" that followed "Địa chỉ: "
# in the "Bên A" address line, leaving just "Địa chỉ: ".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false,
              $true, 1, $false, "", 2)
